$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 36151
$ws.Range("D2").Value = 725
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 32.95
$ws.Range("H2").Value = 43700

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2506
$ws.Range("D3").Value = 2500
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 83.33
$ws.Range("H3").Value = 43503

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1050
$ws.Range("D4").Value = 3031
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 60.62
$ws.Range("H4").Value = 43504

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 2573
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 43543

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 13497
$ws.Range("D6").Value = 648.43
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 64.84
$ws.Range("H6").Value = 43483

$ws.Range("H3:H6").NumberFormat = "m/d/yy"

$ws.Range("E16").Select()
